$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.913.34'
$ws.Range("E2").Value = '  -1.11%  '

$ws.Range("D3").Value = '1.758.83'
$ws.Range("E3").Value = '  -3.68%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.48%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.21%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3752'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.73%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3343'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.67%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.71'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.20%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.119'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.95%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07130'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.27%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.34%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.11'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.47%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.162'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.139'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.14%  '

$ws.Range("D16").Value = '1.757.65'
$ws.Range("E16").Value = '  -3.74%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001048'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.63%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06561'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '80.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.56%  '

$ws.Range("E20").Value = '  +0.44%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.85'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.249'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.33%  '

$ws.Range("D23").Value = '27.925.25'
$ws.Range("E23").Value = '  -1.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.64'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -8.98%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.390'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.55%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.87%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.58%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.311'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -10.65%  '

$ws.Range("B29").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C29").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D29").Value = '1.956.93'
$ws.Range("E29").Value = '  -3.88%  '

$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.274'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -14.96%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '131.17'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.67%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.019'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.753'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.86%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08651'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.35%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.12'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.87%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02328'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6530'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.54%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06163'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.29%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.120'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.75%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2102'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.84%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.209'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.75%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.446'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -10.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.025'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.60%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.46%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.62'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.89%  '

$ws.Range("E46").Value = '  -1.58%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6000'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.84'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.998'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.25%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07168'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.80%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.174'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.56%  '
